$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Id  No." header text to "Id  Number"
$ws.Range("A1").Value = "Id  Number"

# Move the active selection from B8 to B1
$ws.Range("B1").Select()
